# eventbuttons.xlsx — add "WebSocket Command" action rows to the Commands sheet.
#
# Commit message: "adds WebSocket Command actions - adds Probat Sample Roaster setup"
# (the Probat Sample Roaster portion belongs to other files not covered by this
# workbook's diff; only the WebSocket Command rows apply here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# New action block appended after the existing "RC Command" rows (row 88 was
# the previous last row, "range(c,r[,sn])").
#   Row 89: Action="WebSocket Command", Command="send(<json>)", Description=...
#   Row 90: continuation row (no Action), reuses the existing "sleep(<float>)"
#           command/description pair used elsewhere in the sheet (e.g. row 11).
$ws.Range("A89").Value = "WebSocket Command"
$ws.Range("B89").Value = "send(<json>)"
$ws.Range("C89").Value = "If {} substitutions are used, json brackets need to be duplicated to escape them like in send({{ “value”: {}}})"

$ws.Range("B90").Value = "sleep(<float>)"
$ws.Range("C90").Value = "sleep: add a delay of <float> seconds"

# Match the row heights used for the new rows.
$ws.Rows.Item(89).RowHeight = 13.8
$ws.Rows.Item(90).RowHeight = 13.8

# Widen column C slightly to fit the new, longer description text.
$ws.Columns.Item(3).ColumnWidth = 40.4

# Bring the new rows into view and select the last edited cell, mirroring the
# saved cursor position/selection in the source edit.
$ws.Activate()
$ws.Range("C89").Select()
